# Apply updated crypto price/volume figures (cell-by-cell text values).
# Values that look like plain numbers get a leading apostrophe so Excel
# keeps them as text (matching the original inline-string cell type)
# instead of silently converting them to a numeric value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.652.14'
$ws.Range('E2').Value = '  -1.59%  '
$ws.Range('D3').Value = '2.904.46'
$ws.Range('E3').Value = '  -2.42%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''528.13'
$ws.Range('E5').Value = '  -2.83%  '
$ws.Range('D6').Value = '''143.03'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '''0.553'
$ws.Range('E8').Value = '  -3.70%  '
$ws.Range('D9').Value = '2.912.14'
$ws.Range('E9').Value = '  -2.51%  '
$ws.Range('E10').Value = '  -5.30%  '
$ws.Range('E11').Value = '  -2.15%  '
$ws.Range('E12').Value = '  -3.07%  '
$ws.Range('D13').Value = '3.418.40'
$ws.Range('E13').Value = '  -2.27%  '
$ws.Range('D15').Value = '60.634.35'
$ws.Range('E15').Value = '  -1.73%  '
$ws.Range('D16').Value = '''22.83'
$ws.Range('E16').Value = '  -3.81%  '
$ws.Range('D17').Value = '2.913.70'
$ws.Range('E17').Value = '  -2.57%  '
$ws.Range('E18').Value = '  -4.24%  '
$ws.Range('D19').Value = '''5.01'
$ws.Range('E19').Value = '  -3.49%  '
$ws.Range('D20').Value = '''11.68'
$ws.Range('E20').Value = '  -2.87%  '
$ws.Range('D21').Value = '''359.97'
$ws.Range('E21').Value = '  -5.88%  '
$ws.Range('D22').Value = '''6.67'
$ws.Range('E22').Value = '  -0.27%  '
$ws.Range('D24').Value = '''5.68'
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('D25').Value = '''64.86'
$ws.Range('E25').Value = '  -1.34%  '
$ws.Range('E26').Value = '  -3.70%  '
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('D29').Value = '''7.91'
$ws.Range('E29').Value = '  -4.69%  '
$ws.Range('D30').Value = '0.0₃0848'
$ws.Range('E30').Value = '  -10.17%  '
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('D32').Value = '''1.68'
$ws.Range('D33').Value = '''19.79'
$ws.Range('E33').Value = '  -3.45%  '
$ws.Range('D34').Value = '''150.97'
$ws.Range('E34').Value = '  -5.96%  '
$ws.Range('E35').Value = '  -7.48%  '
$ws.Range('D36').Value = '''5.60'
$ws.Range('E36').Value = '  -6.03%  '
$ws.Range('E37').Value = '  -6.83%  '
$ws.Range('E38').Value = '  -5.44%  '
$ws.Range('D39').Value = '''37.91'
$ws.Range('E39').Value = '  +1.36%  '
$ws.Range('E40').Value = '  -4.89%  '
$ws.Range('E41').Value = '  -5.55%  '
$ws.Range('D42').Value = '2.296.78'
$ws.Range('E42').Value = '  -4.73%  '
$ws.Range('E43').Value = '  -2.92%  '
$ws.Range('E44').Value = '  -1.78%  '
$ws.Range('D45').Value = '''20.44'
$ws.Range('E45').Value = '  -7.85%  '
$ws.Range('D46').Value = '''0.997'
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('E47').Value = '  -4.01%  '
$ws.Range('D48').Value = '''0.0238'
$ws.Range('E48').Value = '  -4.04%  '
$ws.Range('D49').Value = '''10.31'
$ws.Range('E49').Value = '  -1.34%  '
$ws.Range('D51').Value = '''249.52'
$ws.Range('E51').Value = '  -7.84%  '
